$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Desired OOXML <col> widths (as seen in Excel's column-width dialog / the
# saved width attribute). The runtime's ColumnWidth setter stores the value
# plus a fixed 5/6-character padding offset (Excel's standard internal
# width-vs-displayed-width quirk), so we subtract that offset here to land
# on the exact target width in the saved file.
$offset = 5/6
$targetWidths = @(5, 15, 20, 13, 12, 15, 10, 18, 22, 22, 10, 14)

for ($i = 0; $i -lt $targetWidths.Length; $i++) {
    $col = $i + 1
    $ws.Columns.Item($col).ColumnWidth = $targetWidths[$i] - $offset
}
